$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the sub_pop_section3 -> sub_pop_section2 values on row 3 (F3, G3)
$ws.Range("F3").Value = "sub_pop_section2"
$ws.Range("G3").Value = "sub_pop_section2_checkbox"

# Update the active selection to F4 as recorded in the saved workbook
$ws.Range("F4").Select()
